$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "88.123.10"
$ws.Range("E2").Value = "  -2.80%  "

$ws.Range("D3").Value = "3.066.75"
$ws.Range("E3").Value = "  -4.26%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "619.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.371"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.23%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.797"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +12.30%  "

$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").Value = "3.063.62"
$ws.Range("E10").Value = "  -4.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.594"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.56%  "

$ws.Range("E12").Value = "  -1.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.69%  "

$ws.Range("D15").Value = "87.845.14"
$ws.Range("E15").Value = "  -2.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.33%  "

$ws.Range("D17").Value = "3.622.78"
$ws.Range("E17").Value = "  -4.53%  "

$ws.Range("D18").Value = "3.067.02"
$ws.Range("E18").Value = "  -4.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000205"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "419.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.73%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "81.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.33%  "

$ws.Range("D28").Value = "3.215.13"
$ws.Range("E28").Value = "  -4.98%  "

$ws.Range("E29").Value = "  +0.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.174"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.86%  "

$ws.Range("E31").Value = "  +8.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "506.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.83%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -12.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.56%  "

$ws.Range("E36").Value = "  -4.55%  "

$ws.Range("E37").Value = "  -7.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.87%  "

$ws.Range("E39").Value = "  +2.15%  "

$ws.Range("E40").Value = "  -0.82%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("E43").Value = "  -3.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "147.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.45%  "

$ws.Range("E45").Value = "  -7.34%  "

$ws.Range("E46").Value = "  +4.16%  "

$ws.Range("E47").Value = "  -3.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0678"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.01%  "

$ws.Range("E49").Value = "  -3.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "158.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.699"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.05%  "
